# chore: update Sheets via scheduled runner
# Refreshes cached market-price / profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the per-job Leve tracking sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2485.1904
$ws.Range("J40").Value = 2874.3333
$ws.Range("L40").Value = 2874.3333
$ws.Range("N40").Value = -3224.3333

$ws.Range("H62").Value = 66668770
$ws.Range("I62").Value = 66668770
$ws.Range("K62").Value = 66668770
$ws.Range("M62").Value = -66668146

$ws.Range("H65").Value = 66668770
$ws.Range("I65").Value = 66668770
$ws.Range("K65").Value = 333343850
$ws.Range("M65").Value = -333340730

$ws.Range("H80").Value = 872.05554
$ws.Range("I80").Value = 700.4
$ws.Range("J80").Value = 1086.625
$ws.Range("K80").Value = 2101.2
$ws.Range("L80").Value = 3259.875
$ws.Range("M80").Value = -1103.2
$ws.Range("N80").Value = -5255.875

$ws.Range("H83").Value = 872.05554
$ws.Range("I83").Value = 700.4
$ws.Range("J83").Value = 1086.625
$ws.Range("K83").Value = 6303.599999999999
$ws.Range("L83").Value = 9779.625
$ws.Range("M83").Value = -1311.599999999999
$ws.Range("N83").Value = -19763.625

$ws.Range("H92").Value = 2043.75
$ws.Range("I92").Value = 2002.2727
$ws.Range("K92").Value = 2002.2727
$ws.Range("M92").Value = -754.2727

$ws.Range("H132").Value = 1421.5883
$ws.Range("I132").Value = 1236.9286
$ws.Range("K132").Value = 3710.7858
$ws.Range("M132").Value = -1180.7858

$ws.Range("H138").Value = 2930.8262
$ws.Range("I138").Value = 1400.2354
$ws.Range("K138").Value = 4200.706200000001
$ws.Range("M138").Value = 939.2937999999995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 153.28572
$ws.Range("I4").Value = 153.28572
$ws.Range("K4").Value = 153.28572
$ws.Range("M4").Value = -37.28572

$ws.Range("H8").Value = 5404
$ws.Range("I8").Value = 6200
$ws.Range("K8").Value = 6200
$ws.Range("M8").Value = -6056

$ws.Range("H61").Value = 4836.5
$ws.Range("I61").Value = 3496.652
$ws.Range("K61").Value = 3496.652
$ws.Range("M61").Value = -3284.652

$ws.Range("H132").Value = 2244.1929
$ws.Range("J132").Value = 4975.4
$ws.Range("L132").Value = 14926.2
$ws.Range("N132").Value = -19986.2

$ws.Range("H136").Value = 4836.5
$ws.Range("I136").Value = 3496.652
$ws.Range("K136").Value = 10489.956
$ws.Range("M136").Value = -7939.956

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 103817.6
$ws.Range("J132").Value = 103817.6
$ws.Range("L132").Value = 103817.6
$ws.Range("N132").Value = -113937.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5219.4
$ws.Range("I16").Value = 4199.4546
$ws.Range("J16").Value = 6466
$ws.Range("K16").Value = 4199.4546
$ws.Range("L16").Value = 6466
$ws.Range("M16").Value = -3912.4546
$ws.Range("N16").Value = -7040

$ws.Range("H22").Value = 726.06665
$ws.Range("I22").Value = 693.5455
$ws.Range("K22").Value = 693.5455
$ws.Range("M22").Value = -343.5454999999999

$ws.Range("H113").Value = 5219.4
$ws.Range("I113").Value = 4199.4546
$ws.Range("J113").Value = 6466
$ws.Range("K113").Value = 4199.4546
$ws.Range("L113").Value = 6466
$ws.Range("M113").Value = -2029.4546
$ws.Range("N113").Value = -10806

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2271.6453
$ws.Range("J2").Value = 4126.353
$ws.Range("L2").Value = 24758.118
$ws.Range("N2").Value = -24984.118

$ws.Range("H15").Value = 1186.8572
$ws.Range("I15").Value = 539.25
$ws.Range("J15").Value = 2050.3333
$ws.Range("K15").Value = 1617.75
$ws.Range("L15").Value = 6150.999899999999
$ws.Range("M15").Value = -1477.75
$ws.Range("N15").Value = -6430.999899999999

$ws.Range("H22").Value = 1045.2
$ws.Range("J22").Value = 1583.6666
$ws.Range("L22").Value = 4750.9998
$ws.Range("N22").Value = -5088.9998

$ws.Range("H27").Value = 1045.2
$ws.Range("J27").Value = 1583.6666
$ws.Range("L27").Value = 4750.9998
$ws.Range("N27").Value = -4954.9998

$ws.Range("H50").Value = 801
$ws.Range("I50").Value = 151.14285
$ws.Range("K50").Value = 453.42855
$ws.Range("M50").Value = 27.57144999999997

$ws.Range("H53").Value = 801
$ws.Range("I53").Value = 151.14285
$ws.Range("K53").Value = 453.42855
$ws.Range("M53").Value = 27.57144999999997

$ws.Range("H69").Value = 1082
$ws.Range("I69").Value = 1095.8334
$ws.Range("J69").Value = 999
$ws.Range("K69").Value = 3287.5002
$ws.Range("L69").Value = 2997
$ws.Range("M69").Value = -2476.5002
$ws.Range("N69").Value = -4619

$ws.Range("H72").Value = 1082
$ws.Range("I72").Value = 1095.8334
$ws.Range("J72").Value = 999
$ws.Range("K72").Value = 9862.5006
$ws.Range("L72").Value = 8991
$ws.Range("M72").Value = -5806.500599999999
$ws.Range("N72").Value = -17103

$ws.Range("H129").Value = 26316926
$ws.Range("I129").Value = 892.4167
$ws.Range("J129").Value = 71430130
$ws.Range("K129").Value = 2677.2501
$ws.Range("L129").Value = 214290390
$ws.Range("M129").Value = 2322.7499
$ws.Range("N129").Value = -214300390

$ws.Range("H132").Value = 1492.25
$ws.Range("I132").Value = 1450.2
$ws.Range("K132").Value = 13051.8
$ws.Range("M132").Value = -10521.8

$ws.Range("H137").Value = 4148.6
$ws.Range("I137").Value = 2117.4
$ws.Range("J137").Value = 6179.8
$ws.Range("K137").Value = 6352.200000000001
$ws.Range("L137").Value = 18539.4
$ws.Range("M137").Value = -1252.200000000001
$ws.Range("N137").Value = -28739.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 1200
$ws.Range("J17").Value = 1200
$ws.Range("L17").Value = 1200
$ws.Range("N17").Value = -1540

$ws.Range("H18").Value = 3501
$ws.Range("I18").Value = 505
$ws.Range("J18").Value = 4999
$ws.Range("K18").Value = 505
$ws.Range("L18").Value = 4999
$ws.Range("M18").Value = -333
$ws.Range("N18").Value = -5343

$ws.Range("H22").Value = 2528.7693
$ws.Range("I22").Value = 1514.9688
$ws.Range("J22").Value = 4150.85
$ws.Range("K22").Value = 1514.9688
$ws.Range("L22").Value = 4150.85
$ws.Range("M22").Value = -1219.9688
$ws.Range("N22").Value = -4740.85

$ws.Range("H27").Value = 2528.7693
$ws.Range("I27").Value = 1514.9688
$ws.Range("J27").Value = 4150.85
$ws.Range("K27").Value = 1514.9688
$ws.Range("L27").Value = 4150.85
$ws.Range("M27").Value = -1407.9688
$ws.Range("N27").Value = -4364.85

$ws.Range("H46").Value = 4458.6055
$ws.Range("I46").Value = 1911.5385
$ws.Range("J46").Value = 5783.08
$ws.Range("K46").Value = 1911.5385
$ws.Range("L46").Value = 5783.08
$ws.Range("M46").Value = -1723.5385
$ws.Range("N46").Value = -6159.08

$ws.Range("H55").Value = 308.05
$ws.Range("I55").Value = 275.64706
$ws.Range("J55").Value = 491.66666
$ws.Range("K55").Value = 275.64706
$ws.Range("L55").Value = 491.66666
$ws.Range("M55").Value = -102.64706
$ws.Range("N55").Value = -837.66666

$ws.Range("H132").Value = 8096.8
$ws.Range("J132").Value = 7493.4443
$ws.Range("L132").Value = 22480.3329
$ws.Range("N132").Value = -27540.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 22197.143
$ws.Range("J62").Value = 4096.75
$ws.Range("L62").Value = 4096.75
$ws.Range("N62").Value = -5344.75

$ws.Range("H65").Value = 22197.143
$ws.Range("J65").Value = 4096.75
$ws.Range("L65").Value = 20483.75
$ws.Range("N65").Value = -26723.75

$ws.Range("H126").Value = 2375.923
$ws.Range("I126").Value = 2375.923
$ws.Range("K126").Value = 7127.768999999999
$ws.Range("M126").Value = -4657.768999999999

$ws.Range("H132").Value = 9758.167
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 9758.167
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 29274.501
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -34334.501

$ws.Range("H136").Value = 4048.2
$ws.Range("I136").Value = 1474.8182
$ws.Range("K136").Value = 4424.4546
$ws.Range("M136").Value = -1874.4546

$ws.Range("H137").Value = 96928.5
$ws.Range("J137").Value = 96928.5
$ws.Range("L137").Value = 96928.5
$ws.Range("N137").Value = -107128.5
